$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New header cells J1:L1 (ketamineTrial, controlTrial, gainTrial), styled like H1
$ws.Range("J1").Value = "ketamineTrial"
$ws.Range("K1").Value = "controlTrial"
$ws.Range("L1").Value = "gainTrial"
$ws.Range("H1").Copy()
$ws.Range("J1:L1").PasteSpecial(-4122)

# Update existing value and add new data row values
$ws.Range("I2").Value = 100
$ws.Range("J2").Value = 100
$ws.Range("K2").Value = 50
$ws.Range("L2").Value = 290

# Update selection to I3 as in the diff
$ws.Range("I3").Select()
